$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 93-97: set column E to "Add redirect"
foreach ($r in 93..97) {
    $ws.Range("E$r").Value = "Add redirect"
}

# Rows 98-99: set column E to "Add redirect when  done"
foreach ($r in 98..99) {
    $ws.Range("E$r").Value = "Add redirect when  done"
}

# Row 100: set D to "Maria" and E to "Add redirect when  done"
$ws.Range("D100").Value = "Maria"
$ws.Range("E100").Value = "Add redirect when  done"

# Row 102: set E to "Add redirect when  done"
$ws.Range("E102").Value = "Add redirect when  done"

# Rows 111-115: set D to "Eugeney" and E to "Add redirect"
foreach ($r in 111..115) {
    $ws.Range("D$r").Value = "Eugeney"
    $ws.Range("E$r").Value = "Add redirect"
}

# Update view: topLeftCell A76, selection E109
$ws.Application.ActiveWindow.ScrollRow = 76
$ws.Range("E109").Select()
